$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow cell edits, then re-protect at the end.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (shared string used by cell A80).
$disclaimerText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-28 for illustrative purposes only and are subject to change."
$ws.Cells.Item(80, 1).Value2 = $disclaimerText

# Refresh Weight (column D) and Percent Change (column E) figures for each holding row (2-77).
$ws.Cells.Item(2, 4).Value2 = 0.07394568140755603
$ws.Cells.Item(2, 5).Value2 = -0.00534802043422733
$ws.Cells.Item(3, 4).Value2 = 0.04539406680870294
$ws.Cells.Item(3, 5).Value2 = -0.002179492339270062
$ws.Cells.Item(4, 4).Value2 = 0.03854021775012462
$ws.Cells.Item(4, 5).Value2 = 0.001484096105250687
$ws.Cells.Item(5, 4).Value2 = 0.03506972764346906
$ws.Cells.Item(5, 5).Value2 = 0.001939522172264896
$ws.Cells.Item(6, 4).Value2 = 0.03320371559098181
$ws.Cells.Item(6, 5).Value2 = -0.002467536864916076
$ws.Cells.Item(7, 4).Value2 = 0.03079571028029109
$ws.Cells.Item(7, 5).Value2 = -0.0006693033160936057
$ws.Cells.Item(8, 4).Value2 = 0.03116999948798736
$ws.Cells.Item(8, 5).Value2 = 0.0004839267196681263
$ws.Cells.Item(9, 4).Value2 = 0.02925906336163519
$ws.Cells.Item(9, 5).Value2 = 0.002606480658728705
$ws.Cells.Item(10, 4).Value2 = 0.02655701382607469
$ws.Cells.Item(10, 5).Value2 = -0.0009631056452807929
$ws.Cells.Item(11, 4).Value2 = 0.02661116630718819
$ws.Cells.Item(11, 5).Value2 = 0.001197031362221823
$ws.Cells.Item(12, 4).Value2 = 0.02348381367817658
$ws.Cells.Item(12, 5).Value2 = -0.00217828418230559
$ws.Cells.Item(13, 4).Value2 = 0.0243175557706144
$ws.Cells.Item(13, 5).Value2 = -0.003760282021151529
$ws.Cells.Item(14, 4).Value2 = 0.02057588165602972
$ws.Cells.Item(14, 5).Value2 = -0.0003779289493577309
$ws.Cells.Item(15, 4).Value2 = 0.02000250244423969
$ws.Cells.Item(15, 5).Value2 = 0.008627715482112164
$ws.Cells.Item(16, 4).Value2 = 0.01847030577273412
$ws.Cells.Item(16, 5).Value2 = 0.01119993507284023
$ws.Cells.Item(17, 4).Value2 = 0.01725730893367441
$ws.Cells.Item(17, 5).Value2 = 0.002399604770978891
$ws.Cells.Item(18, 4).Value2 = 0.0171580449998686
$ws.Cells.Item(18, 5).Value2 = 0.00976042590949433
$ws.Cells.Item(19, 4).Value2 = 0.0171463338231572
$ws.Cells.Item(19, 5).Value2 = -0.0120811419984973
$ws.Cells.Item(20, 4).Value2 = 0.01536206578411551
$ws.Cells.Item(20, 5).Value2 = -0.003244535519125735
$ws.Cells.Item(21, 4).Value2 = 0.01377838677978664
$ws.Cells.Item(21, 5).Value2 = -0.009424436216761922
$ws.Cells.Item(22, 4).Value2 = 0.01525863267140045
$ws.Cells.Item(22, 5).Value2 = 0.007834757834757955
$ws.Cells.Item(23, 4).Value2 = 0.01356622710248296
$ws.Cells.Item(23, 5).Value2 = -0.003922651933701671
$ws.Cells.Item(24, 4).Value2 = 0.01546624841214011
$ws.Cells.Item(24, 5).Value2 = -0.001889992730797152
$ws.Cells.Item(25, 4).Value2 = 0.01343941848105195
$ws.Cells.Item(25, 5).Value2 = 0.0004391881293723809
$ws.Cells.Item(26, 4).Value2 = 0.01100869348754111
$ws.Cells.Item(26, 5).Value2 = 0.002502085070892557
$ws.Cells.Item(27, 4).Value2 = 0.01209905088407907
$ws.Cells.Item(27, 5).Value2 = 0.006427133343658076
$ws.Cells.Item(28, 4).Value2 = 0.01177001366319565
$ws.Cells.Item(28, 5).Value2 = 0.007156048014773875
$ws.Cells.Item(29, 4).Value2 = 0.01200671996688641
$ws.Cells.Item(29, 5).Value2 = 0.007826490681169807
$ws.Cells.Item(30, 4).Value2 = 0.01214964316747231
$ws.Cells.Item(30, 5).Value2 = -0.0008482418260336022
$ws.Cells.Item(31, 4).Value2 = 0.01054694521216413
$ws.Cells.Item(31, 5).Value2 = -0.01056642993244417
$ws.Cells.Item(32, 4).Value2 = 0.01235594725642029
$ws.Cells.Item(32, 5).Value2 = 0
$ws.Cells.Item(33, 4).Value2 = 0.01081124304818695
$ws.Cells.Item(33, 5).Value2 = 0.003553028957186122
$ws.Cells.Item(34, 4).Value2 = 0.01091753368801959
$ws.Cells.Item(34, 5).Value2 = -0.003604253018561931
$ws.Cells.Item(35, 4).Value2 = 0.01081175833996225
$ws.Cells.Item(35, 5).Value2 = 0.004159445407279083
$ws.Cells.Item(36, 4).Value2 = 0.009755035442939423
$ws.Cells.Item(36, 5).Value2 = 0.00829803785979788
$ws.Cells.Item(37, 4).Value2 = 0.01015991424420586
$ws.Cells.Item(37, 5).Value2 = 0.00456463102565885
$ws.Cells.Item(38, 4).Value2 = 0.0088655949940622
$ws.Cells.Item(38, 5).Value2 = -0.008924466988983193
$ws.Cells.Item(39, 4).Value2 = 0.01013906834965957
$ws.Cells.Item(39, 5).Value2 = 0.002069857697283073
$ws.Cells.Item(40, 4).Value2 = 0.009168399179112096
$ws.Cells.Item(40, 5).Value2 = 0.007572080380545687
$ws.Cells.Item(41, 4).Value2 = 0.008763332999018275
$ws.Cells.Item(41, 5).Value2 = 0.009942695860417405
$ws.Cells.Item(42, 4).Value2 = 0.008764457271982568
$ws.Cells.Item(42, 5).Value2 = 0.006734510625561185
$ws.Cells.Item(43, 4).Value2 = 0.009619466861328285
$ws.Cells.Item(43, 5).Value2 = -0.00303874398581927
$ws.Cells.Item(44, 4).Value2 = 0.008955958433567343
$ws.Cells.Item(44, 5).Value2 = -0.0007322788517868384
$ws.Cells.Item(45, 4).Value2 = 0.008853274836161809
$ws.Cells.Item(45, 5).Value2 = -0.006307145275990522
$ws.Cells.Item(46, 4).Value2 = 0.009491112364571366
$ws.Cells.Item(46, 5).Value2 = -0.001658374792703143
$ws.Cells.Item(47, 4).Value2 = 0.008535480344921324
$ws.Cells.Item(47, 5).Value2 = -0.0003951527924129783
$ws.Cells.Item(48, 4).Value2 = 0.008533231798992735
$ws.Cells.Item(48, 5).Value2 = 0.004797979797979579
$ws.Cells.Item(49, 4).Value2 = 0.007976014011064439
$ws.Cells.Item(49, 5).Value2 = 0.003259624702669273
$ws.Cells.Item(50, 4).Value2 = 0.008906303044311017
$ws.Cells.Item(50, 5).Value2 = 0.0002314279102058059
$ws.Cells.Item(51, 4).Value2 = 0.007911883607392825
$ws.Cells.Item(51, 5).Value2 = -0.0008289124668434678
$ws.Cells.Item(52, 4).Value2 = 0.008211221284136146
$ws.Cells.Item(52, 5).Value2 = 0.0006960053854843107
$ws.Cells.Item(53, 4).Value2 = 0.006752664491791847
$ws.Cells.Item(53, 5).Value2 = -0.007700312174817814
$ws.Cells.Item(54, 4).Value2 = 0.007465172482913275
$ws.Cells.Item(54, 5).Value2 = 0.002710843373494143
$ws.Cells.Item(55, 4).Value2 = 0.006716125620452287
$ws.Cells.Item(55, 5).Value2 = -0.003034107553881449
$ws.Cells.Item(56, 4).Value2 = 0.006809580810609239
$ws.Cells.Item(56, 5).Value2 = 0.002700443710659473
$ws.Cells.Item(57, 4).Value2 = 0.0079748428933933
$ws.Cells.Item(57, 5).Value2 = 0.008458646616541277
$ws.Cells.Item(58, 4).Value2 = 0.006532025922549116
$ws.Cells.Item(58, 5).Value2 = 0.002065404475042998
$ws.Cells.Item(59, 4).Value2 = 0.006641080400085651
$ws.Cells.Item(59, 5).Value2 = -0.0006348400203149973
$ws.Cells.Item(60, 4).Value2 = 0.005761524184352847
$ws.Cells.Item(60, 5).Value2 = 0.004813321191622144
$ws.Cells.Item(61, 4).Value2 = 0.005851653400323763
$ws.Cells.Item(61, 5).Value2 = -0.008277562522014725
$ws.Cells.Item(62, 4).Value2 = 0.005860741273451806
$ws.Cells.Item(62, 5).Value2 = 0.002317960195028501
$ws.Cells.Item(63, 4).Value2 = 0.004979030201204106
$ws.Cells.Item(63, 5).Value2 = -0.01678458527773607
$ws.Cells.Item(64, 4).Value2 = 0.004925065098917985
$ws.Cells.Item(64, 5).Value2 = 0.0001521838380764695
$ws.Cells.Item(65, 4).Value2 = 0.004514705466950614
$ws.Cells.Item(65, 5).Value2 = 0.003942890346144257
$ws.Cells.Item(66, 4).Value2 = 0.004469547169551466
$ws.Cells.Item(66, 5).Value2 = -0.0001257703433531487
$ws.Cells.Item(67, 4).Value2 = 0.00453110111434657
$ws.Cells.Item(67, 5).Value2 = -0.02059425593945796
$ws.Cells.Item(68, 4).Value2 = 0.004435116310019956
$ws.Cells.Item(68, 5).Value2 = -0.002471561202826322
$ws.Cells.Item(69, 4).Value2 = 0.004116338079935714
$ws.Cells.Item(69, 5).Value2 = -0.003778222869628589
$ws.Cells.Item(70, 4).Value2 = 0.003530779244365834
$ws.Cells.Item(70, 5).Value2 = 0.002335084646818641
$ws.Cells.Item(71, 4).Value2 = 0.003646438825567597
$ws.Cells.Item(71, 5).Value2 = 0.002119705553628393
$ws.Cells.Item(72, 4).Value2 = 0.002902123278498012
$ws.Cells.Item(72, 5).Value2 = 0.04884426652892571
$ws.Cells.Item(73, 4).Value2 = 0.002360317399121915
$ws.Cells.Item(73, 5).Value2 = -0.002083912197832771
$ws.Cells.Item(74, 4).Value2 = 0.002334178052702075
$ws.Cells.Item(74, 5).Value2 = 0.01264349361804618
$ws.Cells.Item(75, 4).Value2 = 0.001967477687514794
$ws.Cells.Item(75, 5).Value2 = -0.01466666666666672
$ws.Cells.Item(76, 4).Value2 = 0.00200598403654187
$ws.Cells.Item(76, 5).Value2 = -0.005744710662743513
$ws.Cells.Item(77, 4).Value2 = 0.9999999999999999
$ws.Cells.Item(77, 5).Value2 = 0.0001603751961212208

# Restore sheet protection.
$ws.Protect()
